$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 33686531.32457295
$ws.Range("E4").Value = 6872304.585514704
$ws.Range("C7").Value = 525595.3322263303
$ws.Range("G7").Value = 117220857.186105
$ws.Range("I8").Value = 49886784.6835634
$ws.Range("E10").Value = 9447875.952797854
$ws.Range("E14").Value = 163899.5297748352
$ws.Range("E17").Value = 8032350.765089288
$ws.Range("E27").Value = -5273485.495753591
$ws.Range("E29").Value = 9585279.296527933
$ws.Range("C31").Value = 258318.683726668
$ws.Range("G31").Value = 50438406.8982071
$ws.Range("I31").Value = 54015049.9983958
$ws.Range("E33").Value = 6069764.464888245
$ws.Range("C34").Value = 585527.948057324
$ws.Range("G34").Value = 240162116.010487
$ws.Range("E35").Value = -7245394.153791295
$ws.Range("E36").Value = 15674217.66116988
$ws.Range("I36").Value = 140937442.588492
$ws.Range("E37").Value = 19674320.38476177
$ws.Range("E38").Value = 27574170.14620421
$ws.Range("C39").Value = 1986971.395975875
$ws.Range("E39").Value = 9073580.796574228
$ws.Range("G39").Value = 359549616.289326
$ws.Range("E40").Value = 2744269.595138839
$ws.Range("C42").Value = 560592.6138586564
$ws.Range("G42").Value = 138921011.0633
$ws.Range("E43").Value = 8365578.784704204
$ws.Range("E44").Value = 67699177.56819201
$ws.Range("E45").Value = 31097026.64363619
$ws.Range("E46").Value = 37792529.97965544
$ws.Range("E47").Value = 1489812.729699545
$ws.Range("C50").Value = 1318373.742259424
$ws.Range("E50").Value = 99142898.92025533
$ws.Range("G50").Value = 311861775.963436
$ws.Range("E52").Value = 2406525.909258161
$ws.Range("C54").Value = 2557827.666945276
$ws.Range("E54").Value = 24041690.26878396
$ws.Range("G54").Value = 552717606.141346
$ws.Range("E56").Value = 46342348.84790083
$ws.Range("E59").Value = -16038339.70004726
$ws.Range("E60").Value = -7739752.17270211
$ws.Range("B63").Value = 431734.2034979998
$ws.Range("E63").Value = 5398299.655918255
$ws.Range("C64").Value = 4730814.943753129
$ws.Range("G64").Value = 1006967360.25823
$ws.Range("C65").Value = 2436212.450195243
$ws.Range("E65").Value = 16745936.93818601
$ws.Range("G65").Value = 571824545.031608
$ws.Range("E66").Value = -4735994.283433722
$ws.Range("E67").Value = 20982985.68700554
$ws.Range("I67").Value = 303764449.160993
